$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.792.16"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.148.33"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.149.49"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").Value = "3.693.68"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "57.956.22"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "3.154.99"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "358.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.510"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D28").Value = "0.0₃0938"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0672"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +11.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").Value = "3.192.07"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0274"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "2.330.87"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
